# [Feat 2269] Added support of datasets worksheet metadata
# Rename "DATASETS TODO" sheet to "DATASETS", give it a header row with
# the new dataset/parameter metadata columns, size those columns, and
# make DATASETS the active tab (instead of PARAMETERS).

$wb = $excel.ActiveWorkbook

# --- Rename the DATASETS TODO sheet to DATASETS -----------------------
$wsDatasets = $wb.Worksheets.Item("DATASETS TODO")
$wsDatasets.Name = "DATASETS"

# --- Populate header row on the DATASETS sheet -------------------------
$wsDatasets.Cells.Item(1, 1).Value = "ACTION"
$wsDatasets.Cells.Item(1, 2).Value = "TC_OWNER_PATH"
$wsDatasets.Cells.Item(1, 3).Value = "TC_OWNER_ID"
$wsDatasets.Cells.Item(1, 4).Value = "TC_DATASET_ID"
$wsDatasets.Cells.Item(1, 5).Value = "TC_DATASET_NAME"
$wsDatasets.Cells.Item(1, 6).Value = "TC_PARAM_OWNER_ID"
$wsDatasets.Cells.Item(1, 7).Value = "TC_DATASET_PARAM_VALUE"
$wsDatasets.Cells.Item(1, 8).Value = "TC_DATASET_PARAM_NAME"
$wsDatasets.Cells.Item(1, 9).Value = "TC_PARAM_OWNER_PATH"

# --- Size the columns to fit their (bestFit) header text ----------------
$wsDatasets.Columns.Item(1).ColumnWidth = 7.0221354166667
$wsDatasets.Columns.Item(2).ColumnWidth = 16.0221354166667
$wsDatasets.Columns.Item(3).ColumnWidth = 13.0221354166667
$wsDatasets.Columns.Item(4).ColumnWidth = 14.0221354166667
$wsDatasets.Columns.Item(5).ColumnWidth = 17.7369791666667
$wsDatasets.Columns.Item(6).ColumnWidth = 20.8776041666667
$wsDatasets.Columns.Item(7).ColumnWidth = 25.8776041666667
$wsDatasets.Columns.Item(8).ColumnWidth = 25.5924479166667

# --- Record the selection that was active on this sheet -----------------
$wsDatasets.Range("H7").Select()

# --- Make DATASETS the active sheet/tab (was PARAMETERS) ----------------
$wsDatasets.Activate()
